$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.835.60'
$ws.Range('E2').Value = '  +2.46%  '
$ws.Range('D3').Value = '1.869.47'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('D4').Value = '''1.014'
$ws.Range('E4').Value = '  -0.46%  '
$ws.Range('D5').Value = '''313.71'
$ws.Range('E5').Value = '  +0.99%  '
$ws.Range('D6').Value = '''1.012'
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('D7').Value = '''0.4836'
$ws.Range('E7').Value = '  +1.25%  '
$ws.Range('D8').Value = '''0.3819'
$ws.Range('E8').Value = '  +3.07%  '
$ws.Range('D9').Value = '''0.07367'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('D10').Value = '''0.9396'
$ws.Range('E10').Value = '  +0.57%  '
$ws.Range('D11').Value = '''21.00'
$ws.Range('E11').Value = '  +5.14%  '
$ws.Range('D12').Value = '''0.07829'
$ws.Range('E12').Value = '  +0.14%  '
$ws.Range('D13').Value = '1.903.78'
$ws.Range('E13').Value = '  +2.48%  '
$ws.Range('D14').Value = '''5.497'
$ws.Range('D15').Value = '''6.613'
$ws.Range('E15').Value = '  +1.50%  '
$ws.Range('D16').Value = '''91.13'
$ws.Range('E16').Value = '  +1.56%  '
$ws.Range('D17').Value = '''1.015'
$ws.Range('E17').Value = '  -0.48%  '
$ws.Range('D18').Value = '''0.000008862'
$ws.Range('E18').Value = '  +1.71%  '
$ws.Range('E19').Value = '  -0.49%  '
$ws.Range('D20').Value = '27.864.44'
$ws.Range('E20').Value = '  +2.45%  '
$ws.Range('D22').Value = '''5.126'
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').Value = '2.116.97'
$ws.Range('E23').Value = '  +1.29%  '
$ws.Range('E24').Value = '  +1.81%  '
$ws.Range('D25').Value = '''1.950'
$ws.Range('E25').Value = '  +0.46%  '
$ws.Range('D26').Value = '''156.96'
$ws.Range('E26').Value = '  +2.43%  '
$ws.Range('D27').Value = '''18.56'
$ws.Range('E27').Value = '  +0.67%  '
$ws.Range('D28').Value = '''2.053'
$ws.Range('E28').Value = '  +2.86%  '
$ws.Range('D29').Value = '''115.92'
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('D30').Value = '''4.983'
$ws.Range('E30').Value = '  +1.04%  '
$ws.Range('D31').Value = '''0.08919'
$ws.Range('E31').Value = '  +0.52%  '
$ws.Range('D32').Value = '''3.336'
$ws.Range('E32').Value = '  +0.83%  '
$ws.Range('D33').Value = '''1.231'
$ws.Range('E33').Value = '  +4.18%  '
$ws.Range('D34').Value = '''0.7692'
$ws.Range('E34').Value = '  +4.29%  '
$ws.Range('D35').Value = '''4.656'
$ws.Range('E35').Value = '  +2.14%  '
$ws.Range('D36').Value = '''2.739'
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('D37').Value = '''1.135'
$ws.Range('E37').Value = '  +1.48%  '
$ws.Range('D38').Value = '''0.02048'
$ws.Range('E38').Value = '  +2.47%  '
$ws.Range('D39').Value = '''0.5616'
$ws.Range('E39').Value = '  +5.97%  '
$ws.Range('D40').Value = '''0.05363'
$ws.Range('E40').Value = '  +2.21%  '
$ws.Range('D41').Value = '''2.998'
$ws.Range('E41').Value = '  +0.53%  '
$ws.Range('D42').Value = '''7.047'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').Value = '''8.572'
$ws.Range('E43').Value = '  +2.94%  '
$ws.Range('D44').Value = '''0.1533'
$ws.Range('E44').Value = '  +0.17%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D45').Value = '''10.75'
$ws.Range('E45').Value = '  +0.99%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').Value = '''0.4881'
$ws.Range('E46').Value = '  +2.62%  '
$ws.Range('D47').Value = '''105.45'
$ws.Range('E47').Value = '  +3.04%  '
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').Value = '''1.668'
$ws.Range('E49').Value = '  +2.48%  '
$ws.Range('D50').Value = '''68.07'
$ws.Range('E50').Value = '  +3.13%  '
$ws.Range('D51').Value = '''0.06119'
$ws.Range('E51').Value = '  +0.82%  '
